# Shuffled SO 1215 and PredictionSet3
# Multiply the values in H2:K8 by 100 (scale correction)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("H2:K8")
foreach ($cell in $range.Cells) {
    $current = $cell.Value()
    $cell.Value = $current * 100
}

# Update the active selection to reflect the author's last edit location
$ws.Range("I13").Select()
